# Insert a new weekly data row for "Feria Lagunitas de Puerto Montt - Cebollín"
# at row 119, pushing the existing rows 119-225 down to 120-226
# (dimension grows from A1:R225 to A1:R226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 119..225 down by one row, creating an empty row 119.
$ws.Rows.Item(119).Insert()

# Populate the newly inserted row 119 with the new record.
$ws.Cells.Item(119, 1).Value = 4
$ws.Cells.Item(119, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(119, 3).Value = "Los Lagos"
$ws.Cells.Item(119, 4).Value = 44586
$ws.Cells.Item(119, 5).Value = 10
$ws.Cells.Item(119, 6).Value = 100112037
$ws.Cells.Item(119, 7).Value = "Cebollín"
$ws.Cells.Item(119, 8).Value = "Sin especificar"
$ws.Cells.Item(119, 9).Value = "Primera"
$ws.Cells.Item(119, 10).Value = 180
$ws.Cells.Item(119, 11).Value = 6000
$ws.Cells.Item(119, 12).Value = 6000
$ws.Cells.Item(119, 13).Value = 6000
$ws.Cells.Item(119, 14).Value = "`$/paquete 36 unidades"
$ws.Cells.Item(119, 15).Value = "Región Metropolitana"
$ws.Cells.Item(119, 16).Value = 167
$ws.Cells.Item(119, 17).Value = 36
$ws.Cells.Item(119, 18).Value = "Hortaliza"

# Apply the same date style/number format used by the rest of column D
# to the new cell (mirrors the style of D118/D120, style index "2").
$ws.Cells.Item(119, 4).NumberFormat = $ws.Cells.Item(120, 4).NumberFormat
